# Add a "language" column to the manuscript metadata header row.
#
# The new field is inserted immediately before "source_collection" (column
# N), so every column from source_collection through batch shifts one
# position to the right (N:Z -> O:AA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank column at N; this pushes source_collection..batch
# (old N:Z) right by one to O:AA and extends the used range to column AA.
$ws.Columns("N:N").Insert()

# Give the new column the same explicit width as its neighbours
# (~15.33 chars) instead of inheriting the "best fit" auto-size flag from
# the column to its left.
$ws.Columns("N:N").ColumnWidth = 14.5

# Populate the new header cell. Matching the style already used for every
# other header cell in row 1 (t="s" / s="1") happens automatically because
# the inserted column copies formatting from its neighbour.
$ws.Range("N1").Value = "language"

# Reflect the editor's new cursor position/selection.
[void]$ws.Range("K3").Select()
